$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 98
    3  = 30
    4  = 128
    5  = 59
    6  = 72
    7  = 101
    8  = 41
    9  = 26
    10 = 198
    11 = 240
    12 = 47
    13 = 8
    14 = 62
    15 = 20
    16 = 13
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
